$wb = $excel.ActiveWorkbook

# Excel constant values (not predefined in this host, so use literals):
#   xlEdgeLeft   = 7
#   xlEdgeTop    = 8
#   xlEdgeBottom = 9
#   xlEdgeRight  = 10
#   xlContinuous = 1
#   xlNone       = -4142
#   xlThin       = 2

function Set-TopBottomBorder($cell) {
    # Reset to the plain "Normal" style first so the new xf does not
    # inherit the bold font / box border / centered alignment that the
    # cell previously had.
    $cell.Style = "Normal"
    # Top + bottom thin border, no left/right (matches the workbook's
    # existing border id 4). Top is applied before bottom so every
    # intermediate combination (none -> top-only) matches a border
    # permutation that already exists in the workbook (id 2), instead of
    # registering a brand-new, unused border entry.
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(8).Weight = 2
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(9).Weight = 2
    $cell.Borders.Item(7).LineStyle = -4142
    $cell.Borders.Item(10).LineStyle = -4142
}

function Set-TopRightBottomBorder($cell) {
    # Reset to the plain "Normal" style first (see comment above).
    $cell.Style = "Normal"
    # Top + right + bottom thin border, no left (matches the workbook's
    # existing border id 5). Order: top (-> id 2), right (-> id 3),
    # bottom (-> id 5) so every intermediate state already exists in the
    # workbook and no spurious border/style entries get created.
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(8).Weight = 2
    $cell.Borders.Item(10).LineStyle = 1
    $cell.Borders.Item(10).Weight = 2
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(9).Weight = 2
    $cell.Borders.Item(7).LineStyle = -4142
}

# ---- Sheet 1: quality_comparison ----
$ws1 = $wb.Worksheets.Item(1)

Set-TopBottomBorder($ws1.Range("C1"))
Set-TopRightBottomBorder($ws1.Range("D1"))

$ws1.Range("C2").Value = "approach"

# ---- Sheet 2: computational_comparison ----
$ws2 = $wb.Worksheets.Item(2)

Set-TopBottomBorder($ws2.Range("C1"))
Set-TopRightBottomBorder($ws2.Range("D1"))
Set-TopBottomBorder($ws2.Range("F1"))
Set-TopRightBottomBorder($ws2.Range("G1"))

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell at G5
$ws2.Range("G5").ClearContents()
